$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D44").Value = 475.2
$ws1.Range("L49").Value = 380.16
$ws1.Range("D55").Value = "5 de 53"
$ws1.Range("L55").Value = "7 de 53"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F44").Value = 1477.07
$ws2.Range("F49").Value = 380.16
$ws2.Range("F55").Value = 48044.59

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 7197.02
$ws3.Range("E3").Value = 6530.98
$ws3.Range("F3").Value = 0.5242584498834499

$ws3.Range("D15").Value = 12839.3
$ws3.Range("E15").Value = 7850.700000000001
$ws3.Range("F15").Value = 0.6205558240695989

$ws3.Range("D19").Value = 53881.42
$ws3.Range("E19").Value = 51331.45
$ws3.Range("F19").Value = 0.5121181467628437
